# New Measurements added to the Excel file
# Adds two new flow-rate readings to Tabelle1 (row 65 completed + blank
# spacer rows 66-70) and a fresh block of perfusion measurements
# (rows 71-78), then grows Table1 / conditional formatting / dimension to
# cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Cornea Measurements")

# ---------------------------------------------------------------------
# 1. Finish row 65: it already had F/G/J/K values, now also gets H/I and
#    updated start/end timestamps, plus the T-column divisor changes.
# ---------------------------------------------------------------------
$ws.Range("F65").Value = 45931.513888888891
$ws.Range("G65").Value = 45933.553472222222
$ws.Range("H65").Value = 12
$ws.Range("I65").Value = 7

$ws.Range("T65").Formula = "=P65*1440/2937"

# ---------------------------------------------------------------------
# 2. Blank spacer rows 66-70: copy the formatted-but-empty formula cells
#    from row 65 down, then clear their contents so only formatting +
#    (empty) shared-style remains, matching the template rows below the
#    table in the source file.
# ---------------------------------------------------------------------
$ws.Range("F65:G65").Copy()
$ws.Range("F66:G70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N65:R65").Copy()
$ws.Range("N66:R70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("T65:U65").Copy()
$ws.Range("T66:U70").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 66; $r -le 70; $r++) {
    $ws.Range("F$r:G$r").ClearContents()
    $ws.Range("N$r:R$r").ClearContents()
    $ws.Range("T$r:U$r").ClearContents()
}

# ---------------------------------------------------------------------
# 3. New data block, rows 71-78. F/G need the date/time style (s=11)
#    already used on the "Cornea Measurements" sheet (numFmtId 22), so
#    copy formats from there instead of re-typing a NumberFormat string
#    (which would mint a brand-new custom format code).
# ---------------------------------------------------------------------
$ws2.Range("D3").Copy()
$ws.Range("F71:G78").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Pin the N:U formula-column styles (s=3/4/5/1/5/-/1/7, matching row 65)
# onto every new row BEFORE any formula is entered. Formulas that refer
# back to F/G (now date-formatted) would otherwise inherit that date
# format from their precedent cell once evaluated.
$ws.Range("N65:U65").Copy()
$ws.Range("N71:U71").PasteSpecial(-4122)
$ws.Range("N72:U72").PasteSpecial(-4122)
$ws.Range("N73:U73").PasteSpecial(-4122)
$ws.Range("N74:U74").PasteSpecial(-4122)
$ws.Range("N75:U75").PasteSpecial(-4122)
$ws.Range("N76:U76").PasteSpecial(-4122)
$ws.Range("N77:U77").PasteSpecial(-4122)
$ws.Range("N78:U78").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 71 ---
$ws.Range("A71").Value = 1
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 45903.534722222219
$ws.Range("G71").Value = 45905.469444444447
$ws.Range("H71").Value = 11.7
$ws.Range("I71").Value = 7
$ws.Range("J71").Value = 2.5
$ws.Range("N71").Formula = "=G71-F71"
$ws.Range("O71").Formula = "=N71"
$ws.Range("P71").Formula = "=H71-I71"
$ws.Range("Q71").Formula = "=ABS((E71-D71)/0.9982)"
$ws.Range("R71").Formula = "=J71*O71"
$ws.Range("S71").Formula = "=(1-ABS(Q71-R71)/R71)*100"
$ws.Range("T71").Formula = "=P71*1440/2786"
$ws.Range("U71").Formula = "=(1-ABS(T71-J71)/J71)*100%"

# --- Row 72 ---
$ws.Range("A72").Value = 1
$ws.Range("F72").Value = 45905.501388888886
$ws.Range("G72").Value = 45908.460416666669
$ws.Range("H72").Value = 10.1
$ws.Range("I72").Value = 2.7
$ws.Range("J72").Value = 2.5
$ws.Range("N72").Formula = "=G72-F72"
$ws.Range("O72").Formula = "=N72"
$ws.Range("P72").Formula = "=H72-I72"
$ws.Range("Q72").Formula = "=ABS((E72-D72)/0.9982)"
$ws.Range("R72").Formula = "=J72*O72"
$ws.Range("S72").Formula = "=(1-ABS(Q72-R72)/R72)*100"
$ws.Range("T72").Formula = "=P72*1440/4261"
$ws.Range("U72").Formula = "=(1-ABS(T72-J72)/J72)*100%"

# --- Row 73 ---
$ws.Range("A73").Value = 2
$ws.Range("F73").Value = 45909.48541666667
$ws.Range("G73").Value = 45911.393750000003
$ws.Range("H73").Value = 10.7
$ws.Range("I73").Value = 6
$ws.Range("J73").Value = 2.5
$ws.Range("N73").Formula = "=G73-F73"
$ws.Range("O73").Formula = "=N73"
$ws.Range("P73").Formula = "=H73-I73"
$ws.Range("Q73").Formula = "=ABS((E73-D73)/0.9982)"
$ws.Range("R73").Formula = "=J73*O73"
$ws.Range("S73").Formula = "=(1-ABS(Q73-R73)/R73)*100"
$ws.Range("T73").Formula = "=P73*1440/2748"
$ws.Range("U73").Formula = "=(1-ABS(T73-J73)/J73)*100%"

# --- Row 74 ---
$ws.Range("A74").Value = 3
$ws.Range("F74").Value = 45912.517361111109
$ws.Range("G74").Value = 45915.452777777777
$ws.Range("H74").Value = 11.2
$ws.Range("I74").Value = 2.4
$ws.Range("J74").Value = 2.5
$ws.Range("N74").Formula = "=G74-F74"
$ws.Range("O74").Formula = "=N74"
$ws.Range("P74").Formula = "=H74-I74"
$ws.Range("Q74").Formula = "=ABS((E74-D74)/0.9982)"
$ws.Range("R74").Formula = "=J74*O74"
$ws.Range("S74").Formula = "=(1-ABS(Q74-R74)/R74)*100"
$ws.Range("T74").Formula = "=P74*1440/4227"
$ws.Range("U74").Formula = "=(1-ABS(T74-J74)/J74)*100%"

# --- Row 75 ---
$ws.Range("A75").Value = 3
$ws.Range("F75").Value = 45915.701388888891
$ws.Range("G75").Value = 45916.731944444444
$ws.Range("H75").Value = 11.3
$ws.Range("I75").Value = 8.3000000000000007
$ws.Range("J75").Value = 2.5
$ws.Range("N75").Formula = "=G75-F75"
$ws.Range("O75").Formula = "=N75"
$ws.Range("P75").Formula = "=H75-I75"
$ws.Range("Q75").Formula = "=ABS((E75-D75)/0.9982)"
$ws.Range("R75").Formula = "=J75*O75"
$ws.Range("S75").Formula = "=(1-ABS(Q75-R75)/R75)*100"
$ws.Range("T75").Formula = "=P75*1440/1484"
$ws.Range("U75").Formula = "=(1-ABS(T75-J75)/J75)*100%"

# --- Row 76 ---
$ws.Range("A76").Value = 3
$ws.Range("F76").Value = 45919.415972222225
$ws.Range("G76").Value = 45920.629166666666
$ws.Range("H76").Value = 11.6
$ws.Range("I76").Value = 7.5
$ws.Range("J76").Value = 2.5
$ws.Range("N76").Formula = "=G76-F76"
$ws.Range("O76").Formula = "=N76"
$ws.Range("P76").Formula = "=H76-I76"
$ws.Range("Q76").Formula = "=ABS((E76-D76)/0.9982)"
$ws.Range("R76").Formula = "=J76*O76"
$ws.Range("S76").Formula = "=(1-ABS(Q76-R76)/R76)*100"
$ws.Range("T76").Formula = "=P76*1440/1747"
$ws.Range("U76").Formula = "=(1-ABS(T76-J76)/J76)*100%"

# --- Row 77 ---
$ws.Range("A77").Value = 4
$ws.Range("F77").Value = 45919.415972222225
$ws.Range("G77").Value = 45922.427777777775
$ws.Range("H77").Value = 11.6
$ws.Range("I77").Value = 2.2000000000000002
$ws.Range("J77").Value = 2.5
$ws.Range("N77").Formula = "=G77-F77"
$ws.Range("O77").Formula = "=N77"
$ws.Range("P77").Formula = "=H77-I77"
$ws.Range("Q77").Formula = "=ABS((E77-D77)/0.9982)"
$ws.Range("R77").Formula = "=J77*O77"
$ws.Range("S77").Formula = "=(1-ABS(Q77-R77)/R77)*100"
$ws.Range("T77").Formula = "=P77*1440/4337"
$ws.Range("U77").Formula = "=(1-ABS(T77-J77)/J77)*100%"

# --- Row 78 ---
$ws.Range("A78").Value = 5
$ws.Range("F78").Value = 45930.557638888888
$ws.Range("G78").Value = 45930.695138888892
$ws.Range("H78").Value = 9.75
$ws.Range("I78").Value = 9.5
$ws.Range("J78").Value = 2.5
$ws.Range("N78").Formula = "=G78-F78"
$ws.Range("O78").Formula = "=N78"
$ws.Range("P78").Formula = "=H78-I78"
$ws.Range("Q78").Formula = "=ABS((E78-D78)/0.9982)"
$ws.Range("R78").Formula = "=J78*O78"
$ws.Range("S78").Formula = "=(1-ABS(Q78-R78)/R78)*100"
$ws.Range("T78").Formula = "=P78*1440/198"
$ws.Range("U78").Formula = "=(1-ABS(T78-J78)/J78)*100%"

# ---------------------------------------------------------------------
# 4. Grow Table1 to cover the new rows, and the two color-scale
#    conditional formats that track column O and column U.
# ---------------------------------------------------------------------
$t1 = $ws.ListObjects.Item("Table1")
$t1.Resize($ws.Range("A3:K76"))

$ocf = $ws.Range("O4:O65").FormatConditions.Item(1)
$ocf.ModifyAppliesToRange($ws.Range("O4:O78"))

$ucf = $ws.Range("U4:U65").FormatConditions.Item(1)
$ucf.ModifyAppliesToRange($ws.Range("U4:U78"))

# ---------------------------------------------------------------------
# 5. View cosmetics: zoom + active selection on both sheets.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 80
$ws.Range("I66").Select()

$ws2.Activate()
$ws2.Range("B22").Select()

$ws.Activate()
